$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.056197333333333
$ws.Range("H2").Value = 3.168592
$ws.Range("I2").Value = 0.01247237710445079
$ws.Range("J2").Value = 0.01398563433468744
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 128.780808988496
$ws.Range("R2").Value = 1159.027280896464
$ws.Range("S2").Value = 0.002846486003352171
$ws.Range("T2").Value = 0.003384132514460595
$ws.Range("G3").Value = 1.056197333333333
$ws.Range("H3").Value = 3.168592
$ws.Range("I3").Value = 0.01247237710445079
$ws.Range("J3").Value = 0.01398563433468744
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 156.2267209077866
$ws.Range("R3").Value = 1406.04048817008
$ws.Range("S3").Value = 0.003453132325433252
$ws.Range("T3").Value = 0.00410536267013845
$ws.Range("G4").Value = 1.056197333333333
$ws.Range("H4").Value = 3.168592
$ws.Range("I4").Value = 0.01247237710445079
$ws.Range("J4").Value = 0.01398563433468744
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 88.19772592994843
$ws.Range("R4").Value = 793.7795333695359
$ws.Range("S4").Value = 0.001949464321267897
$ws.Range("T4").Value = 0.002317680672806501
$ws.Range("G5").Value = 1.056197333333333
$ws.Range("H5").Value = 3.168592
$ws.Range("I5").Value = 0.01247237710445079
$ws.Range("J5").Value = 0.01398563433468744
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 96.18626724326667
$ws.Range("R5").Value = 577.1176034596
$ws.Range("S5").Value = 0.002126037765821986
$ws.Range("T5").Value = 0.001685070298797975
$ws.Range("G6").Value = 1.056197333333333
$ws.Range("H6").Value = 3.168592
$ws.Range("I6").Value = 0.01247237710445079
$ws.Range("J6").Value = 0.01398563433468744
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 94.88415284431998
$ws.Range("R6").Value = 853.9573755988799
$ws.Range("S6").Value = 0.002097256688575483
$ws.Range("T6").Value = 0.002493388178483923
$ws.Range("I7").Value = 0.6620593097549599
$ws.Range("J7").Value = 0.7423861014276285
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 6835.948977054551
$ws.Range("R7").Value = 61523.54079349095
$ws.Range("S7").Value = 0.1510973042928593
$ws.Range("T7").Value = 0.1796366817552023
$ws.Range("I8").Value = 0.6620593097549599
$ws.Range("J8").Value = 0.7423861014276285
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 8292.834168121833
$ws.Range("R8").Value = 74635.5075130965
$ws.Range("S8").Value = 0.1832993329758328
$ws.Range("T8").Value = 0.2179210548978446
$ws.Range("I9").Value = 0.6620593097549599
$ws.Range("J9").Value = 0.7423861014276285
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 4681.715847919756
$ws.Range("R9").Value = 42135.4426312778
$ws.Range("S9").Value = 0.1034815570537849
$ws.Range("T9").Value = 0.1230272347941688
$ws.Range("I10").Value = 0.6620593097549599
$ws.Range("J10").Value = 0.7423861014276285
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 5105.763974716459
$ws.Range("R10").Value = 30634.58384829876
$ws.Range("S10").Value = 0.1128544369662131
$ws.Range("T10").Value = 0.08944698108211156
$ws.Range("I11").Value = 0.6620593097549599
$ws.Range("J11").Value = 0.7423861014276285
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 5036.6450767735
$ws.Range("R11").Value = 45329.8056909615
$ws.Range("S11").Value = 0.1113266784662697
$ws.Range("T11").Value = 0.1323541488983012
$ws.Range("G12").Value = 0.04559766666666667
$ws.Range("H12").Value = 0.136793
$ws.Range("I12").Value = 0.0005384517417354892
$ws.Range("J12").Value = 0.0006037813885615125
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 5.559665998009
$ws.Range("R12").Value = 50.03699398208099
$ws.Range("S12").Value = 0.0001228871877024728
$ws.Range("T12").Value = 0.0001460982161952717
$ws.Range("G13").Value = 0.04559766666666667
$ws.Range("H13").Value = 0.136793
$ws.Range("I13").Value = 0.0005384517417354892
$ws.Range("J13").Value = 0.0006037813885615125
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 6.744548314563333
$ws.Range("R13").Value = 60.70093483106999
$ws.Range("S13").Value = 0.0001490770443758587
$ws.Range("T13").Value = 0.000177234833558959
$ws.Range("G14").Value = 0.04559766666666667
$ws.Range("H14").Value = 0.136793
$ws.Range("I14").Value = 0.0005384517417354892
$ws.Range("J14").Value = 0.0006037813885615125
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 3.807631756671556
$ws.Range("R14").Value = 34.268685810044
$ws.Range("S14").Value = 0.00008416137921802476
$ws.Range("T14").Value = 0.000100057846600389
$ws.Range("G15").Value = 0.04559766666666667
$ws.Range("H15").Value = 0.136793
$ws.Range("I15").Value = 0.0005384517417354892
$ws.Range("J15").Value = 0.0006037813885615125
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 4.152509396920834
$ws.Range("R15").Value = 24.915056381525
$ws.Range("S15").Value = 0.00009178432695029433
$ws.Range("T15").Value = 0.0000727470817901047
$ws.Range("G16").Value = 0.04559766666666667
$ws.Range("H16").Value = 0.136793
$ws.Range("I16").Value = 0.0005384517417354892
$ws.Range("J16").Value = 0.0006037813885615125
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 4.09629511153
$ws.Range("R16").Value = 36.86665600377
$ws.Range("S16").Value = 0.00009054180348883861
$ws.Range("T16").Value = 0.0001076434104167881
$ws.Range("G17").Value = 27.488287
$ws.Range("H17").Value = 54.976574
$ws.Range("I17").Value = 0.3246024872429512
$ws.Range("J17").Value = 0.2426573888143015
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 3351.612171180093
$ws.Range("R17").Value = 20109.67302708056
$ws.Range("S17").Value = 0.07408182328456371
$ws.Range("T17").Value = 0.05871630415245919
$ws.Range("G18").Value = 27.488287
$ws.Range("H18").Value = 54.976574
$ws.Range("I18").Value = 0.3246024872429512
$ws.Range("J18").Value = 0.2426573888143015
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 4065.911554452709
$ws.Range("R18").Value = 24395.46932671626
$ws.Range("S18").Value = 0.08987022539710379
$ws.Range("T18").Value = 0.07122998941855062
$ws.Range("G19").Value = 27.488287
$ws.Range("H19").Value = 54.976574
$ws.Range("I19").Value = 0.3246024872429512
$ws.Range("J19").Value = 0.2426573888143015
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 2295.408562960865
$ws.Range("R19").Value = 13772.45137776519
$ws.Range("S19").Value = 0.05073619584907635
$ws.Range("T19").Value = 0.04021285890291851
$ws.Range("G20").Value = 27.488287
$ws.Range("H20").Value = 54.976574
$ws.Range("I20").Value = 0.3246024872429512
$ws.Range("J20").Value = 0.2426573888143015
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 2503.316033848738
$ws.Range("R20").Value = 10013.26413539495
$ws.Range("S20").Value = 0.055331645361492
$ws.Range("T20").Value = 0.02923676887938522
$ws.Range("G21").Value = 27.488287
$ws.Range("H21").Value = 54.976574
$ws.Range("I21").Value = 0.3246024872429512
$ws.Range("J21").Value = 0.2426573888143015
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 2469.42758026581
$ws.Range("R21").Value = 14816.56548159486
$ws.Range("S21").Value = 0.05458259735071525
$ws.Range("T21").Value = 0.04326146746098791
$ws.Range("G22").Value = 0.027723
$ws.Range("H22").Value = 0.08316900000000001
$ws.Range("I22").Value = 0.0003273741559027063
$ws.Range("J22").Value = 0.0003670940348210247
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 3.380230431297
$ws.Range("R22").Value = 30.422073881673
$ws.Range("S22").Value = 0.000074714382417426
$ws.Range("T22").Value = 0.00008882649362719254
$ws.Range("G23").Value = 0.027723
$ws.Range("H23").Value = 0.08316900000000001
$ws.Range("I23").Value = 0.0003273741559027063
$ws.Range("J23").Value = 0.0003670940348210247
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 4.10062897059
$ws.Range("R23").Value = 36.90566073531
$ws.Range("S23").Value = 0.00009063759624904632
$ws.Range("T23").Value = 0.0001077572965887514
$ws.Range("G24").Value = 0.027723
$ws.Range("H24").Value = 0.08316900000000001
$ws.Range("I24").Value = 0.0003273741559027063
$ws.Range("J24").Value = 0.0003670940348210247
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 2.315008264828
$ws.Range("R24").Value = 20.835074383452
$ws.Range("S24").Value = 0.0000511694147228579
$ws.Range("T24").Value = 0.00006083433394916224
$ws.Range("G25").Value = 0.027723
$ws.Range("H25").Value = 0.08316900000000001
$ws.Range("I25").Value = 0.0003273741559027063
$ws.Range("J25").Value = 0.0003670940348210247
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 2.5246909858875
$ws.Range("R25").Value = 15.148145915325
$ws.Range("S25").Value = 0.00005580410319335806
$ws.Range("T25").Value = 0.00004422961734446366
$ws.Range("G26").Value = 0.027723
$ws.Range("H26").Value = 0.08316900000000001
$ws.Range("I26").Value = 0.0003273741559027063
$ws.Range("J26").Value = 0.0003670940348210247
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 2.49051317049
$ws.Range("R26").Value = 22.41461853441
$ws.Range("S26").Value = 0.00005504865932001797
$ws.Range("T26").Value = 0.00006544629331145489

Write-Host "Applied 280 cell updates"